$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (A2:H8) contents but keep formatting, and remove stale hyperlinks
$ws.Range("A2:H8").ClearContents()
$ws.Hyperlinks.Delete()

# Column width updates (COM ColumnWidth reports ~0.8333 below the stored
# XML "character width"; subtract that offset so the saved width lands on
# the exact target integer)
$ws.Columns.Item(2).ColumnWidth = 50.166666666666664
$ws.Columns.Item(8).ColumnWidth = 17.166666666666668

$data = @(
    @("2026-01-12 18:29:46", "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5423720", 385, "🔥AI,Ai ◆効率化"),
    @("2026-01-12 18:29:46", "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5434128", 368, "🔥AI,Ai ◆開発"),
    @("2026-01-12 18:29:46", "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5427956", 310, "🔥AI,Ai"),
    @("2026-01-12 18:29:46", "Shopee APIを使用した「商品動画の一括紐付けツール」の開発依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469483", 308, "🔥API ◆ツール,開発"),
    @("2026-01-12 18:29:46", "【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5439158", 303, "🔥AI,Ai"),
    @("2026-01-12 18:29:46", "【Zapier設定のみ!作業時間~1時間】スプレッドシート・Gドライブ自動化構築(設計済)", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469379", 255, "🔥API ◆自動化"),
    @("2026-01-12 18:29:46", "【募集】Python / Docker 日次データ スクレイピングシステム構築", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469627", 248, "🔥Python ◆スクレイピング"),
    @("2026-01-12 18:29:46", "【フルリモート可】Webアプリ開発経験者募集!経営管理システムの開発", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469430", 150, "◆開発 ◇アプリ"),
    @("2026-01-12 18:29:46", "【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5454210", 108, "◆開発 ◇アプリ"),
    @("2026-01-12 18:29:46", "【未経験相談可能】JavaまたはJavascriptエンジニアを募集!", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469522", 85, "★Java"),
    @("2026-01-12 18:29:46", "フロント実装済み!音楽権利マーケットプレイス「HITOON」のバックエンド・決済機能実装", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469298", 18, $null),
    @("2026-01-12 18:29:46", "【急募】Microsoft Accessで物流納品先別仕分けリスト作成", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469531", 10, $null),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[5])
    $ws.Cells.Item($r, 6).Style = $ws.Cells.Item(2, 6).Style
    $ws.Cells.Item($r, 7).Value = [double]$row[6]
    if ($row[7] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row[7]
    }
}

$ws.Range("A1").Select()
